# Allowing resource optimized OpenStack
# - Rename the "LW OpenStack" box (slide 1) to just "OpenStack" so that the
#   box can represent any resource optimized OpenStack deployment, not only
#   the "LW" (light-weight) variant.
# - Refresh the cached text of every "datetimeFigureOut" date field
#   (slide master, all slide layouts and the notes master) so that it shows
#   the date of this edit instead of the stale one.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. "LW OpenStack" -> "OpenStack" (slide 1, inside the big background group)
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $top = $s1.Shapes.Item($i)
    if ($top.Type -eq 6) {
        # msoGroup -- walk its GroupItems looking for the run of text
        for ($k = 1; $k -le $top.GroupItems.Count; $k++) {
            $sh = $top.GroupItems.Item($k)
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.TextRange.Text -eq "LW OpenStack") {
                    $sh.TextFrame.TextRange.Text = "OpenStack"
                }
            }
        }
    } elseif ($top.HasTextFrame) {
        if ($top.TextFrame.TextRange.Text -eq "LW OpenStack") {
            $top.TextFrame.TextRange.Text = "OpenStack"
        }
    }
}

# ---------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" field text: 2018/8/29 -> 2018/10/18
# ---------------------------------------------------------------------
$oldDate = "2018/8/29"
$newDate = "2018/10/18"

function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout that hangs off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
